$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from H1 onto the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-25 for columns I (I0) and J (IF)
$data = @{
    2  = @(8, 9)
    3  = @(7, 7)
    4  = @(6, 6)
    5  = @(8, 9)
    6  = @(7, 8)
    7  = @(1, 1)
    8  = @(7, 8)
    9  = @(5, 6)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(7, 8)
    14 = @(9, 9)
    15 = @(5, 5)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(8, 8)
    19 = @(8, 8)
    20 = @(5, 6)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(3, 4)
    24 = @(4, 4)
    25 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
